# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G (header "K") previously held the pitcher's per-game strikeout-ish
# "Strike#" stat; it is being regenerated to hold the true strikeout count (K)
# for each start. Replace the existing numeric values in G2:G42 with the
# newly computed K values, row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 6
    4  = 9
    5  = 7
    6  = 5
    7  = 4
    8  = 5
    9  = 9
    10 = 11
    11 = 8
    12 = 5
    13 = 5
    14 = 6
    15 = 8
    16 = 7
    17 = 10
    18 = 7
    19 = 5
    20 = 8
    21 = 8
    22 = 3
    23 = 7
    24 = 3
    25 = 9
    26 = 5
    27 = 2
    28 = 6
    29 = 6
    30 = 7
    31 = 6
    32 = 8
    33 = 4
    34 = 5
    35 = 6
    36 = 3
    37 = 4
    38 = 2
    39 = 5
    40 = 3
    41 = 1
    42 = 5
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
